# Atualizado por script em 02-01-2024 20:45
# Applies two changes to the only worksheet:
#   1) Rows 47 and 48 had their match data swapped (index/date/opponent-
#      independent columns stay put, the rest trade places).
#   2) 11 new match rows (93-103) are appended after the previous last
#      row (92), matching the columns: Indice, pais, torneio, temporada,
#      data_partida, home, home_ft_gols, away, away_ft_gols,
#      home_opening_odds, home_opening_data_hora, home_closing_odds,
#      home_closing_data_hora, draw_opening_odds, draw_opening_data_hora,
#      draw_closing_odds, draw_closing_data_hora, away_opening_odds,
#      away_opening_data_hora, away_closing_odds, away_closing_data_hora,
#      url_partida.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Swap the match-specific columns between row 47 and row 48.
#    Columns A (Indice), B-E (pais/torneio/temporada/data_partida) and
#    K/O/S (the *_opening_data_hora columns) are identical between the
#    two rows already and are left untouched.
# ---------------------------------------------------------------------
$swapCols = @("F", "H", "I", "J", "L", "M", "N", "P", "Q", "R", "T", "U", "V")

foreach ($col in $swapCols) {
    $a = $ws.Range($col + "47").Value()
    $b = $ws.Range($col + "48").Value()
    $ws.Range($col + "47").Value = $b
    $ws.Range($col + "48").Value = $a
}

# ---------------------------------------------------------------------
# 2) Append 11 new rows of match data (rows 93-103).
# ---------------------------------------------------------------------
$newRows = @(
    @(92, "morocco", "botola-pro", "2023-2024", 45283.75, "Wydad", 1, "Youssoufia Berrechid", 0, 1.24, "11/11/2023 08:43", 1.25, "23/12/2023 17:59", 4.97, "11/11/2023 08:43", 5.22, "23/12/2023 17:59", 9.77, "11/11/2023 08:43", 13.18, "23/12/2023 17:59", "https://www.betexplorer.com/football/morocco/botola-pro/wydad-athletic-youssoufia-berrechid/MsNTL6p2/"),
    @(93, "morocco", "botola-pro", "2023-2024", 45284.75, "Berkane", 1, "Raja Casablanca", 1, 2.62, "23/12/2023 04:41", 2.79, "24/12/2023 17:59", 2.69, "23/12/2023 04:41", 2.66, "24/12/2023 17:59", 2.75, "23/12/2023 04:41", 2.96, "24/12/2023 17:59", "https://www.betexplorer.com/football/morocco/botola-pro/berkane-raja-casablanca/hSExYk6U/"),
    @(94, "morocco", "botola-pro", "2023-2024", 45286.75, "Hassania Agadir", 0, "Wydad", 0, 5.6, "24/12/2023 18:12", 5.79, "26/12/2023 17:58", 3.38, "24/12/2023 18:12", 2.88, "26/12/2023 17:58", 1.58, "24/12/2023 18:12", 1.81, "26/12/2023 17:58", "https://www.betexplorer.com/football/morocco/botola-pro/hassania-agadir-wydad-athletic/OrR832sa/"),
    @(95, "morocco", "botola-pro", "2023-2024", 45288.66666666666, "FUS Rabat", 1, "Youssoufia Berrechid", 1, 1.6, "27/12/2023 04:12", 1.44, "28/12/2023 15:52", 3.34, "27/12/2023 04:12", 3.91, "28/12/2023 15:52", 5.53, "27/12/2023 04:12", 8.51, "28/12/2023 15:52", "https://www.betexplorer.com/football/morocco/botola-pro/fus-rabat-youssoufia-berrechid/QyRxAzJ9/"),
    @(96, "morocco", "botola-pro", "2023-2024", 45288.66666666666, "Raja Casablanca", 6, "IR Tanger", 1, 1.53, "27/12/2023 04:12", 1.36, "28/12/2023 15:34", 3.6, "27/12/2023 04:12", 4.4, "28/12/2023 15:50", 5.67, "27/12/2023 04:12", 9.4, "28/12/2023 15:50", "https://www.betexplorer.com/football/morocco/botola-pro/raja-casablanca-ir-tanger/4SPt9fYF/"),
    @(97, "morocco", "botola-pro", "2023-2024", 45288.75, "Chabab Mohammedia", 1, "Berkane", 1, 4.25, "27/12/2023 06:12", 7.16, "28/12/2023 17:57", 3, "27/12/2023 06:12", 3.22, "28/12/2023 17:58", 1.85, "27/12/2023 06:12", 1.61, "28/12/2023 17:57", "https://www.betexplorer.com/football/morocco/botola-pro/chabab-mohammedia-berkane/Y1Up8EmM/"),
    @(98, "morocco", "botola-pro", "2023-2024", 45288.83333333334, "Moghreb Tetouan", 1, "Union Touarga", 1, 2.12, "27/12/2023 08:12", 2.43, "28/12/2023 19:53", 2.83, "27/12/2023 08:12", 2.79, "28/12/2023 19:41", 3.56, "27/12/2023 08:12", 3.32, "28/12/2023 19:53", "https://www.betexplorer.com/football/morocco/botola-pro/moghreb-tetouan-union-touarga/ChTl7Y2S/"),
    @(99, "morocco", "botola-pro", "2023-2024", 45289.66666666666, "Renaissance Zemamra", 0, "Mouloudia Oujda", 1, 1.93, "28/12/2023 02:12", 1.81, "29/12/2023 15:17", 2.99, "28/12/2023 02:12", 3.17, "29/12/2023 15:17", 3.95, "28/12/2023 02:12", 4.09, "29/12/2023 15:16", "https://www.betexplorer.com/football/morocco/botola-pro/renaissance-zemamra-mouloudia-oujda/IZ5H2hBq/"),
    @(100, "morocco", "botola-pro", "2023-2024", 45289.75, "Hassania Agadir", 2, "Maghreb Fez", 1, 2.85, "28/12/2023 02:12", 2.55, "29/12/2023 17:59", 2.71, "28/12/2023 02:12", 2.71, "29/12/2023 17:59", 2.59, "28/12/2023 02:12", 3.22, "29/12/2023 17:57", "https://www.betexplorer.com/football/morocco/botola-pro/hassania-agadir-maghreb-fez/U99L1CQk/"),
    @(101, "morocco", "botola-pro", "2023-2024", 45289.75, "Olympique de Safi", 2, "Wydad", 1, 3.76, "28/12/2023 02:12", 4.17, "29/12/2023 17:57", 2.91, "28/12/2023 02:12", 2.89, "29/12/2023 17:51", 2.01, "28/12/2023 02:12", 2.06, "29/12/2023 17:57", "https://www.betexplorer.com/football/morocco/botola-pro/olympique-de-safi-wydad-athletic/OW8P0Wud/"),
    @(102, "morocco", "botola-pro", "2023-2024", 45289.83333333334, "Jeunesse Sportive Soualem", 2, "FAR Rabat", 2, 5.79, "28/12/2023 02:12", 7.5, "29/12/2023 19:56", 3.61, "28/12/2023 02:12", 4.12, "29/12/2023 19:56", 1.52, "28/12/2023 02:12", 1.44, "29/12/2023 19:56", "https://www.betexplorer.com/football/morocco/botola-pro/jeunesse-sportive-soualem-far-rabat/2N7Taje2/")
)

$lastExistingRow = 92
$destRow = $lastExistingRow + 1

foreach ($rowVals in $newRows) {
    # Clone row 92's per-column formatting (Indice column is bordered/
    # bold, data_partida column carries the datetime number format) onto
    # the new row before writing values into it.
    $srcRange = $ws.Range("A" + $lastExistingRow + ":V" + $lastExistingRow)
    $dstRange = $ws.Range("A" + $destRow + ":V" + $destRow)
    $srcRange.Copy($dstRange)

    $col = 1
    foreach ($v in $rowVals) {
        $ws.Cells.Item($destRow, $col).Value = $v
        $col = $col + 1
    }

    $destRow = $destRow + 1
}
